$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values such as "243.42" or
# "1.000" are not silently reinterpreted as numbers, matching the original
# inline-string/text storage used throughout the sheet.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Row 16 / 17 swap: WrappedBTC and Dai traded places in the ranking ---
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '27.909.40'
$ws.Range("E16").Value = '  +5.53%  '

$ws.Range("B17").Value = 'Dai'
$ws.Range("C17").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D17").Value = '0.9996'
$ws.Range("E17").Value = '  +0.12%  '

# --- Remaining price / volume(1h) refreshes ---
$ws.Range("D2").Value = '27.945.13'

$ws.Range("D3").Value = '1.779.41'
$ws.Range("E3").Value = '  +3.43%  '

$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("D5").Value = '243.42'
$ws.Range("E5").Value = '  +1.10%  '

$ws.Range("E6").Value = '  +0.20%  '

$ws.Range("D7").Value = '0.4880'
$ws.Range("E7").Value = '  -1.03%  '

$ws.Range("D8").Value = '0.2664'
$ws.Range("E8").Value = '  +2.32%  '

$ws.Range("D9").Value = '0.06244'
$ws.Range("E9").Value = '  +0.55%  '

$ws.Range("D10").Value = '1.780.36'
$ws.Range("E10").Value = '  +3.48%  '

$ws.Range("D11").Value = '16.32'
$ws.Range("E11").Value = '  +3.39%  '

$ws.Range("D12").Value = '0.07001'
$ws.Range("E12").Value = '  -0.04%  '

$ws.Range("D13").Value = '0.6238'
$ws.Range("E13").Value = '  +2.62%  '

$ws.Range("D14").Value = '4.620'
$ws.Range("E14").Value = '  +2.87%  '

$ws.Range("D15").Value = '79.60'
$ws.Range("E15").Value = '  +3.55%  '

$ws.Range("D18").Value = '0.9992'
$ws.Range("E18").Value = '  +0.14%  '

$ws.Range("D19").Value = '0.000007201'
$ws.Range("E19").Value = '  +0.58%  '

$ws.Range("D20").Value = '11.84'
$ws.Range("E20").Value = '  +4.21%  '

$ws.Range("D21").Value = '2.009.45'
$ws.Range("E21").Value = '  +3.23%  '

$ws.Range("D22").Value = '4.577'
$ws.Range("E22").Value = '  +3.62%  '

$ws.Range("D23").Value = '8.670'
$ws.Range("E23").Value = '  +1.79%  '

$ws.Range("D24").Value = '5.212'
$ws.Range("E24").Value = '  +2.40%  '

$ws.Range("D25").Value = '141.89'
$ws.Range("E25").Value = '  +2.98%  '

$ws.Range("D26").Value = '15.60'
$ws.Range("E26").Value = '  +1.92%  '

$ws.Range("E27").Value = '  +7.16%  '

$ws.Range("D28").Value = '108.88'
$ws.Range("E28").Value = '  +2.86%  '

$ws.Range("D29").Value = '1.387'
$ws.Range("E29").Value = '  -1.09%  '

$ws.Range("D30").Value = '4.178'
$ws.Range("E30").Value = '  +6.51%  '

$ws.Range("D31").Value = '0.08231'
$ws.Range("E31").Value = '  +3.41%  '

$ws.Range("D32").Value = '3.787'
$ws.Range("E32").Value = '  +3.75%  '

$ws.Range("D33").Value = '0.04766'
$ws.Range("E33").Value = '  +6.04%  '

$ws.Range("D34").Value = '1.071'
$ws.Range("E34").Value = '  +7.09%  '

$ws.Range("D35").Value = '2.602'
$ws.Range("E35").Value = '  -0.34%  '

$ws.Range("D36").Value = '0.6429'
$ws.Range("E36").Value = '  +3.21%  '

$ws.Range("D37").Value = '0.9434'
$ws.Range("E37").Value = '  +0.73%  '

$ws.Range("D38").Value = '2.584'
$ws.Range("E38").Value = '  +7.08%  '

$ws.Range("D39").Value = '2.049'
$ws.Range("E39").Value = '  +2.71%  '

$ws.Range("D40").Value = '5.910'
$ws.Range("E40").Value = '  +7.11%  '

$ws.Range("E41").Value = '  +1.92%  '

$ws.Range("D42").Value = '1.000'
$ws.Range("E42").Value = '  +0.23%  '

$ws.Range("D43").Value = '100.02'
$ws.Range("E43").Value = '  +0.72%  '

$ws.Range("D44").Value = '0.3963'
$ws.Range("E44").Value = '  +3.37%  '

$ws.Range("D45").Value = '7.224'
$ws.Range("E45").Value = '  +4.54%  '

$ws.Range("D46").Value = '0.1195'
$ws.Range("E46").Value = '  +3.34%  '

$ws.Range("D47").Value = '0.05414'

$ws.Range("D48").Value = '8.005'

$ws.Range("D49").Value = '1.287'
$ws.Range("E49").Value = '  +5.01%  '

$ws.Range("D50").Value = '30.50'
$ws.Range("E50").Value = '  +1.50%  '

$ws.Range("D51").Value = '52.63'
$ws.Range("E51").Value = '  +2.27%  '

